$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at M:N, shifting the existing M (Bemerkung) and
# N (Gegner) columns to O and P.
$ws.Range("M1:N1").EntireColumn.Insert()

# Headers for the two new columns.
$ws.Range("M1").Value = "VPIP"
$ws.Range("N1").Value = "Alkohol"

# Data entered for the first three data rows only.
$ws.Range("N2").Value = "3 Bier"
$ws.Range("N3").Value = "Nichts"
$ws.Range("M4").Value = 0.51
$ws.Range("M4").Style = "Percent"
$ws.Range("N4").Value = "1 Bier"

# Restore the column width for the two newly inserted columns to match
# the column immediately to their left (Excel inherits this on insert).
$ws.Range("M1:N1").EntireColumn.ColumnWidth = 10.166667

$ws.Range("F1").Select()
